$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 761.1111
$ws.Range("I43").Value = 620
$ws.Range("J43").Value = 937.5
$ws.Range("K43").Value = 620
$ws.Range("L43").Value = 937.5
$ws.Range("M43").Value = -551
$ws.Range("N43").Value = -1075.5

# Row 49: Going Nowhere Fast | Paralyzing Potion
$ws.Range("H49").Value = 1185
$ws.Range("I49").Value = 917
$ws.Range("J49").Value = 1252
$ws.Range("K49").Value = 2751
$ws.Range("L49").Value = 3756
$ws.Range("M49").Value = -2615
$ws.Range("N49").Value = -4028

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 2815.0715
$ws.Range("I116").Value = 2240.5
$ws.Range("J116").Value = 4251.5
$ws.Range("K116").Value = 2240.5
$ws.Range("L116").Value = 4251.5
$ws.Range("M116").Value = 1201.5
$ws.Range("N116").Value = -11135.5

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 1647.1786
$ws.Range("I129").Value = 1262
$ws.Range("J129").Value = 1693.4
$ws.Range("K129").Value = 3786
$ws.Range("L129").Value = 5080.200000000001
$ws.Range("M129").Value = 1214
$ws.Range("N129").Value = -15080.2

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 620.0345
$ws.Range("I2").Value = 606.7037
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 606.7037
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -493.7037
$ws.Range("N2").Value = -1026

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 8783.559999999999
$ws.Range("I32").Value = 6436.375
$ws.Range("J32").Value = 25996.25
$ws.Range("K32").Value = 6436.375
$ws.Range("L32").Value = 25996.25
$ws.Range("M32").Value = -6149.375
$ws.Range("N32").Value = -26570.25

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1218.425
$ws.Range("I45").Value = 1343.4286
$ws.Range("J45").Value = 1151.1154
$ws.Range("K45").Value = 1343.4286
$ws.Range("L45").Value = 1151.1154
$ws.Range("M45").Value = -966.4286
$ws.Range("N45").Value = -1905.1154

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 620.0345
$ws.Range("I116").Value = 606.7037
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 606.7037
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1687.2963
$ws.Range("N116").Value = -5388

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 620.0345
$ws.Range("I3").Value = 606.7037
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 606.7037
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -492.7037
$ws.Range("N3").Value = -1028

# Row 100: And My Axe | Doman Iron War Axe
$ws.Range("H100").Value = 19900
$ws.Range("J100").Value = 19900
$ws.Range("L100").Value = 19900
$ws.Range("N100").Value = -22064

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 1818.8182
$ws.Range("I105").Value = 1464.6154
$ws.Range("K105").Value = 1464.6154
$ws.Range("M105").Value = 282.3846000000001

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 1196.5625
$ws.Range("I107").Value = 1075.6364
$ws.Range("J107").Value = 1462.6
$ws.Range("K107").Value = 1075.6364
$ws.Range("L107").Value = 1462.6
$ws.Range("M107").Value = 844.3635999999999
$ws.Range("N107").Value = -5302.6

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 6314.3125
$ws.Range("I7").Value = 8375.083000000001
$ws.Range("K7").Value = 8375.083000000001
$ws.Range("M7").Value = -8262.083000000001

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 1597.079
$ws.Range("I31").Value = 971.918
$ws.Range("J31").Value = 4139.4
$ws.Range("K31").Value = 971.918
$ws.Range("L31").Value = 4139.4
$ws.Range("M31").Value = -676.918
$ws.Range("N31").Value = -4729.4

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 1597.079
$ws.Range("I34").Value = 971.918
$ws.Range("J34").Value = 4139.4
$ws.Range("K34").Value = 971.918
$ws.Range("L34").Value = 4139.4
$ws.Range("M34").Value = -769.918
$ws.Range("N34").Value = -4543.4

# Row 64: Almost as Fun as Slingshotting Birds | Cedar Longbow
$ws.Range("H64").Value = 29933.334
$ws.Range("J64").Value = 29933.334
$ws.Range("L64").Value = 29933.334
$ws.Range("N64").Value = -30429.334

# Row 67: Living Bow to Mouth (L) | Cedar Longbow
$ws.Range("H67").Value = 29933.334
$ws.Range("J67").Value = 29933.334
$ws.Range("L67").Value = 29933.334
$ws.Range("N67").Value = -31649.334

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up | Kukuru Butter
$ws.Range("H12").Value = 120.695656
$ws.Range("J12").Value = 126.13636
$ws.Range("L12").Value = 378.40908
$ws.Range("N12").Value = -724.40908

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1029.409
$ws.Range("I131").Value = 551.1111
$ws.Range("J131").Value = 1152.4
$ws.Range("K131").Value = 1653.3333
$ws.Range("L131").Value = 3457.2
$ws.Range("M131").Value = 3386.6667
$ws.Range("N131").Value = -13537.2

$ws = $wb.Worksheets.Item("GSM")
# Row 12: Horn of Plenty | Bone Armillae
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("K12").Value = 400
$ws.Range("M12").Value = -260

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 2459.818
$ws.Range("I113").Value = 1998
$ws.Range("J113").Value = 2723.7144
$ws.Range("K113").Value = 1998
$ws.Range("L113").Value = 2723.7144
$ws.Range("M113").Value = 172
$ws.Range("N113").Value = -7063.7144

$ws = $wb.Worksheets.Item("LTW")
# Row 11: A Thorn in One's Hide | Leather Mitts
$ws.Range("H11").Value = 80007
$ws.Range("J11").Value = 80007
$ws.Range("L11").Value = 80007
$ws.Range("N11").Value = -80287

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 1245
$ws.Range("I22").Value = 490
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 490
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -195
$ws.Range("N22").Value = -2590

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 1245
$ws.Range("I27").Value = 490
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 490
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -383
$ws.Range("N27").Value = -2214

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 3631.7144
$ws.Range("I40").Value = 3185.6
$ws.Range("J40").Value = 4747
$ws.Range("K40").Value = 3185.6
$ws.Range("L40").Value = 4747
$ws.Range("M40").Value = -3049.6
$ws.Range("N40").Value = -5019

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 24000
$ws.Range("I122").Value = 27500
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 82500
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -80050
$ws.Range("N122").Value = -34900

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 1583
$ws.Range("I136").Value = 1387.2609
$ws.Range("K136").Value = 4161.7827
$ws.Range("M136").Value = -1611.7827

$ws = $wb.Worksheets.Item("WVR")
# Row 17: Making Gloves Out of Nothing at All | Hempen Bracers
$ws.Range("H17").Value = 1149
$ws.Range("I17").Value = 1149
$ws.Range("K17").Value = 1149
$ws.Range("M17").Value = -977

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 2635.7856
$ws.Range("I81").Value = 2755.6667
$ws.Range("J81").Value = 2420
$ws.Range("K81").Value = 5511.3334
$ws.Range("L81").Value = 4840
$ws.Range("M81").Value = -4450.3334
$ws.Range("N81").Value = -6962

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 2635.7856
$ws.Range("I84").Value = 2755.6667
$ws.Range("J84").Value = 2420
$ws.Range("K84").Value = 27556.667
$ws.Range("L84").Value = 24200
$ws.Range("M84").Value = -22252.667
$ws.Range("N84").Value = -34808

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 4136
$ws.Range("I122").Value = 2520.6
$ws.Range("J122").Value = 6828.3335
$ws.Range("K122").Value = 7561.799999999999
$ws.Range("L122").Value = 20485.0005
$ws.Range("M122").Value = -5111.799999999999
$ws.Range("N122").Value = -25385.0005

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 2962.3333
$ws.Range("I126").Value = 2625.625
$ws.Range("K126").Value = 7876.875
$ws.Range("M126").Value = -5406.875

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2100.4595
$ws.Range("I132").Value = 1944
$ws.Range("J132").Value = 2305.8125
$ws.Range("K132").Value = 5832
$ws.Range("L132").Value = 6917.4375
$ws.Range("M132").Value = -3302
$ws.Range("N132").Value = -11977.4375
